$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 17: One for the Road
$ws.Range("H17").Value = 377378.38
$ws.Range("J17").Value = 377378.38
$ws.Range("L17").Value = 1132135.14
$ws.Range("N17").Value = -1132471.14

# ALC row 40: Stuck in the Moment
$ws.Range("H40").Value = 1288.1052
$ws.Range("I40").Value = 908.125
$ws.Range("J40").Value = 1564.4546
$ws.Range("K40").Value = 908.125
$ws.Range("L40").Value = 1564.4546
$ws.Range("M40").Value = -733.125
$ws.Range("N40").Value = -1914.4546

# ALC row 76: Warding Off Temptation
$ws.Range("H76").Value = 24393510
$ws.Range("I76").Value = 27030140
$ws.Range("K76").Value = 27030140
$ws.Range("M76").Value = -27029825

# ALC row 79: The Garden of Arcane Delights (L)
$ws.Range("H79").Value = 24393510
$ws.Range("I79").Value = 27030140
$ws.Range("K79").Value = 27030140
$ws.Range("M79").Value = -27029048

# ALC row 112: Making Ends Meet
$ws.Range("H112").Value = 1617.5
$ws.Range("J112").Value = 1683.2354
$ws.Range("L112").Value = 5049.706200000001
$ws.Range("N112").Value = -7265.706200000001

# ALC row 135: For Tired Minds
$ws.Range("H135").Value = 13515001
$ws.Range("I135").Value = 637.9231
$ws.Range("J135").Value = 20835280
$ws.Range("K135").Value = 5741.3079
$ws.Range("L135").Value = 187517520
$ws.Range("M135").Value = -3206.3079
$ws.Range("N135").Value = -187522590

# ALC row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 22772.87
$ws.Range("I137").Value = 33197.805
$ws.Range("J137").Value = 1228
$ws.Range("K137").Value = 99593.41500000001
$ws.Range("L137").Value = 3684
$ws.Range("M137").Value = -97043.41500000001
$ws.Range("N137").Value = -8784

# ALC row 138: All-night Crafting
$ws.Range("H138").Value = 1290.38
$ws.Range("I138").Value = 650.9056399999999
$ws.Range("J138").Value = 2011.4894
$ws.Range("K138").Value = 1952.71692
$ws.Range("L138").Value = 6034.468199999999
$ws.Range("M138").Value = 3187.28308
$ws.Range("N138").Value = -16314.4682

# ALC row 141: Remedy for Reason
$ws.Range("H141").Value = 1884.2142
$ws.Range("I141").Value = 1106.4736
$ws.Range("K141").Value = 3319.4208
$ws.Range("M141").Value = 1860.5792

$ws = $wb.Worksheets.Item("ARM")
# ARM row 5: The Alloyed Truth
$ws.Range("H5").Value = 10101009
$ws.Range("I5").Value = 10101009
$ws.Range("K5").Value = 10101009
$ws.Range("M5").Value = -10100897

# ARM row 32: Ingot We Trust
$ws.Range("H32").Value = 17958.988
$ws.Range("I32").Value = 19776.254
$ws.Range("K32").Value = 19776.254
$ws.Range("M32").Value = -19489.254

# ARM row 45: Hollow Hallmarks
$ws.Range("H45").Value = 876
$ws.Range("I45").Value = 832.2
$ws.Range("J45").Value = 1095
$ws.Range("K45").Value = 832.2
$ws.Range("L45").Value = 1095
$ws.Range("M45").Value = -455.2
$ws.Range("N45").Value = -1849

# ARM row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 1833.3334
$ws.Range("I102").Value = 1500
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = 122
$ws.Range("N102").Value = -5744

# ARM row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2440119.8
$ws.Range("I132").Value = 3191591.5
$ws.Range("J132").Value = 722470.4399999999
$ws.Range("K132").Value = 9574774.5
$ws.Range("L132").Value = 2167411.32
$ws.Range("M132").Value = -9572244.5
$ws.Range("N132").Value = -2172471.32

$ws = $wb.Worksheets.Item("BSM")
# BSM row 4: Mending Fences
$ws.Range("H4").Value = 10101009
$ws.Range("I4").Value = 10101009
$ws.Range("K4").Value = 10101009
$ws.Range("M4").Value = -10100894

# BSM row 61: I Maul Right
$ws.Range("H61").Value = 44981.4
$ws.Range("J61").Value = 44981.4
$ws.Range("L61").Value = 44981.4
$ws.Range("N61").Value = -45607.4

# BSM row 105: Ingot to Wing It
$ws.Range("H105").Value = 724680
$ws.Range("I105").Value = 1138038.4
$ws.Range("J105").Value = 1302.75
$ws.Range("K105").Value = 1138038.4
$ws.Range("L105").Value = 1302.75
$ws.Range("M105").Value = -1136291.4
$ws.Range("N105").Value = -4796.75

$ws = $wb.Worksheets.Item("CRP")
# CRP row 4: A Clogful of Camaraderie
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = $null

# CRP row 7: Gridania's Got Talent
$ws.Range("H7").Value = 854.5714
$ws.Range("I7").Value = 1102.4
$ws.Range("J7").Value = 235
$ws.Range("K7").Value = 1102.4
$ws.Range("L7").Value = 235
$ws.Range("M7").Value = -989.4000000000001
$ws.Range("N7").Value = -461

# CRP row 31: Wall Not Found
$ws.Range("H31").Value = 174321.61
$ws.Range("I31").Value = 205769.94
$ws.Range("J31").Value = 17080
$ws.Range("K31").Value = 205769.94
$ws.Range("L31").Value = 17080
$ws.Range("M31").Value = -205474.94
$ws.Range("N31").Value = -17670

# CRP row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 174321.61
$ws.Range("I34").Value = 205769.94
$ws.Range("J34").Value = 17080
$ws.Range("K34").Value = 205769.94
$ws.Range("L34").Value = 17080
$ws.Range("M34").Value = -205567.94
$ws.Range("N34").Value = -17484

# CRP row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1003.5
$ws.Range("I132").Value = 1074.8636
$ws.Range("J132").Value = 807.25
$ws.Range("K132").Value = 3224.5908
$ws.Range("L132").Value = 2421.75
$ws.Range("M132").Value = -694.5907999999999
$ws.Range("N132").Value = -7481.75

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4: In Hot Water
$ws.Range("H4").Value = 673993.7
$ws.Range("I4").Value = 2020304
$ws.Range("J4").Value = 838.5
$ws.Range("K4").Value = 6060912
$ws.Range("L4").Value = 2515.5
$ws.Range("M4").Value = -6060800
$ws.Range("N4").Value = -2739.5

# CUL row 5: What a Sap
$ws.Range("H5").Value = 3082.05
$ws.Range("I5").Value = 549.88464
$ws.Range("J5").Value = 7784.643
$ws.Range("K5").Value = 1649.65392
$ws.Range("L5").Value = 23353.929
$ws.Range("M5").Value = -1537.65392
$ws.Range("N5").Value = -23577.929

# CUL row 122: Salt of the North
$ws.Range("H122").Value = 409.5238
$ws.Range("I122").Value = 312.5
$ws.Range("J122").Value = 432.35294
$ws.Range("K122").Value = 2812.5
$ws.Range("L122").Value = 3891.17646
$ws.Range("M122").Value = -362.5
$ws.Range("N122").Value = -8791.176459999999

# CUL row 125: At Any Temperature
$ws.Range("H125").Value = 7800
$ws.Range("J125").Value = 7800
$ws.Range("L125").Value = 23400
$ws.Range("N125").Value = -33240

# CUL row 131: The Mountain Steeped
$ws.Range("H131").Value = 26596572
$ws.Range("J131").Value = 32052190
$ws.Range("L131").Value = 96156570
$ws.Range("N131").Value = -96166650

# CUL row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 3082.05
$ws.Range("I135").Value = 549.88464
$ws.Range("J135").Value = 7784.643
$ws.Range("K135").Value = 4948.96176
$ws.Range("L135").Value = 70061.787
$ws.Range("M135").Value = -2413.96176
$ws.Range("N135").Value = -75131.787

$ws = $wb.Worksheets.Item("GSM")
# GSM row 2: Copper and Robbers
$ws.Range("H2").Value = 1683522.5
$ws.Range("I2").Value = 2525258.8
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 2525258.8
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = -2525145.8
$ws.Range("N2").Value = -276

# GSM row 5: Hora at Me
$ws.Range("H5").Value = 10574.75
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 11942.571
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 11942.571
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -12166.571

# GSM row 12: Horn of Plenty
$ws.Range("H12").Value = 9595.875
$ws.Range("I12").Value = 24666.666
$ws.Range("J12").Value = 6118
$ws.Range("K12").Value = 24666.666
$ws.Range("L12").Value = 6118
$ws.Range("M12").Value = -24526.666
$ws.Range("N12").Value = -6398

$ws = $wb.Worksheets.Item("LTW")
# LTW row 2: Red in the Head
$ws.Range("H2").Value = 70002
$ws.Range("J2").Value = 70002
$ws.Range("L2").Value = 70002
$ws.Range("N2").Value = -70226

# LTW row 40: Best Served Toad
$ws.Range("H40").Value = 129116.25
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 129116.25
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 129116.25
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = -129388.25

# LTW row 136: Respect for Br'aax
$ws.Range("H136").Value = 324499.16
$ws.Range("I136").Value = 527826.5
$ws.Range("J136").Value = 2564.1667
$ws.Range("K136").Value = 1583479.5
$ws.Range("L136").Value = 7692.500100000001
$ws.Range("M136").Value = -1580929.5
$ws.Range("N136").Value = -12792.5001

$ws = $wb.Worksheets.Item("WVR")
# WVR row 17: Making Gloves Out of Nothing at All
$ws.Range("H17").Value = 2500626
$ws.Range("I17").Value = 3333501.2
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 3333501.2
$ws.Range("L17").Value = 2000
$ws.Range("M17").Value = -3333329.2
$ws.Range("N17").Value = -2344

# WVR row 100: Of Great Import
$ws.Range("H100").Value = 993.1429000000001
$ws.Range("I100").Value = 327
$ws.Range("J100").Value = 1104.1666
$ws.Range("K100").Value = 654
$ws.Range("L100").Value = 2208.3332
$ws.Range("M100").Value = -113
$ws.Range("N100").Value = -3290.3332

# WVR row 126: A Polished Purchase
$ws.Range("H126").Value = 668.7826
$ws.Range("I126").Value = 532.4761999999999
$ws.Range("K126").Value = 1597.4286
$ws.Range("M126").Value = 872.5714000000003

# WVR row 132: Comfy Cabins
$ws.Range("H132").Value = 3600.2
$ws.Range("I132").Value = 529.45
$ws.Range("J132").Value = 7694.533
$ws.Range("K132").Value = 1588.35
$ws.Range("L132").Value = 23083.599
$ws.Range("M132").Value = 941.6499999999999
$ws.Range("N132").Value = -28143.599

# WVR row 136: Weaving the Envelope
$ws.Range("H136").Value = 1851800.5
$ws.Range("I136").Value = 2381574.2
$ws.Range("K136").Value = 7144722.600000001
$ws.Range("M136").Value = -7142172.600000001
